# Weekly price update: a new week's price record for "Espinaca" (Femacal de
# La Calera) is inserted at the top of the data block (row 170), pushing the
# existing rows 170-194 down to 171-195.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 170 - this shifts rows 170:194 down to 171:195
# and extends the used range to A1:R195.
$ws.Rows.Item(170).Insert()

# Populate the new row 170 with the new weekly record. Columns A, B, C, E, F,
# G, H, I, N, Q, R are identical on every data row of this table, so reuse
# those constants; D, J, K, L, M, O, P carry the new week's figures.
$ws.Cells.Item(170, 1).Value = 3
$ws.Cells.Item(170, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(170, 3).Value = "Coquimbo"
$ws.Cells.Item(170, 4).Value = 44474
$ws.Cells.Item(170, 5).Value = 5
$ws.Cells.Item(170, 6).Value = 100112012
$ws.Cells.Item(170, 7).Value = "Espinaca"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 230
$ws.Cells.Item(170, 11).Value = 2500
$ws.Cells.Item(170, 12).Value = 2800
$ws.Cells.Item(170, 13).Value = 2643
$ws.Cells.Item(170, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(170, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(170, 16).Value = 881
$ws.Cells.Item(170, 17).Value = 3
$ws.Cells.Item(170, 18).Value = "Hortaliza"
